$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)
$newws = $wb.Worksheets.Add($null, $sheet1)
$newws.Name = "Выходная мощность -7 дБм"

$newws.Cells.Item(1,1).Value = "OSNR, dB"
$newws.Cells.Item(1,2).Value = "BER"
$newws.Cells.Item(1,3).Value = "Q-factor"
$newws.Cells.Item(1,6).Value = "OSNR, dB"
$newws.Cells.Item(1,7).Value = "BER"
$newws.Cells.Item(1,8).Value = "Q-factor"
$newws.Cells.Item(1,11).Value = "Sequence Length, bit"
$newws.Cells.Item(1,12).Value = 1024
$newws.Cells.Item(2,1).Value = 10
$newws.Cells.Item(2,2).Value = 1
$newws.Cells.Item(2,3).Value = 0
$newws.Cells.Item(2,6).Value = 13
$newws.Cells.Item(2,7).Value = 1
$newws.Cells.Item(2,8).Value = 0
$newws.Cells.Item(2,11).Value = "APD M"
$newws.Cells.Item(2,12).Value = 10
$newws.Cells.Item(3,1).Value = 11
$newws.Cells.Item(3,2).Value = 1
$newws.Cells.Item(3,3).Value = 0
$newws.Cells.Item(3,6).Value = 13.05
$newws.Cells.Item(3,7).Value = 1
$newws.Cells.Item(3,8).Value = 0
$newws.Cells.Item(3,11).Value = "TIA Gain, Om"
$newws.Cells.Item(3,12).Value = 600
$newws.Cells.Item(4,1).Value = 12
$newws.Cells.Item(4,2).Value = 1
$newws.Cells.Item(4,3).Value = 0
$newws.Cells.Item(4,6).Value = 13.1
$newws.Cells.Item(4,7).Value = 1
$newws.Cells.Item(4,8).Value = 0
$newws.Cells.Item(4,11).Value = "LIA Vpp, V"
$newws.Cells.Item(4,12).Value = 0.5
$newws.Cells.Item(5,1).Value = 13
$newws.Cells.Item(5,2).Value = 1
$newws.Cells.Item(5,3).Value = 0
$newws.Cells.Item(5,6).Value = 13.15
$newws.Cells.Item(5,7).Value = 1
$newws.Cells.Item(5,8).Value = 0
$newws.Cells.Item(6,1).Value = 14
$newws.Cells.Item(6,2).Value = 1
$newws.Cells.Item(6,3).Value = 0
$newws.Cells.Item(6,6).Value = 13.2
$newws.Cells.Item(6,7).Value = 1
$newws.Cells.Item(6,8).Value = 0
$newws.Cells.Item(7,1).Value = 15
$newws.Cells.Item(7,2).Value = 1
$newws.Cells.Item(7,3).Value = 0
$newws.Cells.Item(7,6).Value = 13.25
$newws.Cells.Item(7,7).Value = 1
$newws.Cells.Item(7,8).Value = 0
$newws.Cells.Item(8,1).Value = 16
$newws.Cells.Item(8,2).Value = 0
$newws.Cells.Item(8,3).Value = 40.74942828
$newws.Cells.Item(8,6).Value = 13.3
$newws.Cells.Item(8,7).Value = 1
$newws.Cells.Item(8,8).Value = 0
$newws.Cells.Item(9,1).Value = 17
$newws.Cells.Item(9,2).Value = 0
$newws.Cells.Item(9,3).Value = 42.99213023
$newws.Cells.Item(9,6).Value = 13.35
$newws.Cells.Item(9,7).Value = 1
$newws.Cells.Item(9,8).Value = 0
$newws.Cells.Item(10,1).Value = 18
$newws.Cells.Item(10,2).Value = 0
$newws.Cells.Item(10,3).Value = 64.52858408
$newws.Cells.Item(10,6).Value = 13.4
$newws.Cells.Item(10,7).Value = 1
$newws.Cells.Item(10,8).Value = 0
$newws.Cells.Item(11,1).Value = 19
$newws.Cells.Item(11,2).Value = 0
$newws.Cells.Item(11,3).Value = 67.48624979
$newws.Cells.Item(11,6).Value = 13.45
$newws.Cells.Item(11,7).Value = 1
$newws.Cells.Item(11,8).Value = 0
$newws.Cells.Item(12,1).Value = 20
$newws.Cells.Item(12,2).Value = 0
$newws.Cells.Item(12,3).Value = 74.13930622
$newws.Cells.Item(12,6).Value = 13.5
$newws.Cells.Item(12,7).Value = 1
$newws.Cells.Item(12,8).Value = 0
$newws.Cells.Item(13,1).Value = 21
$newws.Cells.Item(13,2).Value = 0
$newws.Cells.Item(13,3).Value = 76.27230078
$newws.Cells.Item(13,6).Value = 13.55
$newws.Cells.Item(13,7).Value = 1
$newws.Cells.Item(13,8).Value = 0
$newws.Cells.Item(14,1).Value = 22
$newws.Cells.Item(14,2).Value = 0
$newws.Cells.Item(14,3).Value = 77.92226651
$newws.Cells.Item(14,6).Value = 13.6
$newws.Cells.Item(14,7).Value = 1
$newws.Cells.Item(14,8).Value = 0
$newws.Cells.Item(15,1).Value = 23
$newws.Cells.Item(15,2).Value = 0
$newws.Cells.Item(15,3).Value = 79.11697408
$newws.Cells.Item(15,6).Value = 13.65
$newws.Cells.Item(15,7).Value = 1
$newws.Cells.Item(15,8).Value = 0
$newws.Cells.Item(16,1).Value = 24
$newws.Cells.Item(16,2).Value = 0
$newws.Cells.Item(16,3).Value = 79.38331063
$newws.Cells.Item(16,6).Value = 13.7
$newws.Cells.Item(16,7).Value = 1
$newws.Cells.Item(16,8).Value = 0
$newws.Cells.Item(17,1).Value = 25
$newws.Cells.Item(17,2).Value = 0
$newws.Cells.Item(17,3).Value = 79.4340575
$newws.Cells.Item(17,6).Value = 13.75
$newws.Cells.Item(17,7).Value = 1
$newws.Cells.Item(17,8).Value = 0
$newws.Cells.Item(18,1).Value = 26
$newws.Cells.Item(18,2).Value = 0
$newws.Cells.Item(18,3).Value = 79.57094632
$newws.Cells.Item(18,6).Value = 13.8
$newws.Cells.Item(18,7).Value = 1
$newws.Cells.Item(18,8).Value = 0
$newws.Cells.Item(19,1).Value = 27
$newws.Cells.Item(19,2).Value = 0
$newws.Cells.Item(19,3).Value = 80.11623534
$newws.Cells.Item(19,6).Value = 13.85
$newws.Cells.Item(19,7).Value = 1
$newws.Cells.Item(19,8).Value = 0
$newws.Cells.Item(20,1).Value = 28
$newws.Cells.Item(20,2).Value = 0
$newws.Cells.Item(20,3).Value = 79.84530655
$newws.Cells.Item(20,6).Value = 13.9
$newws.Cells.Item(20,7).Value = 1
$newws.Cells.Item(20,8).Value = 0
$newws.Cells.Item(21,1).Value = 29
$newws.Cells.Item(21,2).Value = 0
$newws.Cells.Item(21,3).Value = 80.32141827
$newws.Cells.Item(21,6).Value = 13.95
$newws.Cells.Item(21,7).Value = 1
$newws.Cells.Item(21,8).Value = 0
$newws.Cells.Item(22,1).Value = 30
$newws.Cells.Item(22,2).Value = 0
$newws.Cells.Item(22,3).Value = 79.67266047
$newws.Cells.Item(22,6).Value = 14
$newws.Cells.Item(22,7).Value = 1
$newws.Cells.Item(22,8).Value = 0
$newws.Cells.Item(23,6).Value = 14.05
$newws.Cells.Item(23,7).Value = 1
$newws.Cells.Item(23,8).Value = 0
$newws.Cells.Item(24,6).Value = 14.1
$newws.Cells.Item(24,7).Value = 1
$newws.Cells.Item(24,8).Value = 0
$newws.Cells.Item(25,6).Value = 14.15
$newws.Cells.Item(25,7).Value = 1
$newws.Cells.Item(25,8).Value = 0
$newws.Cells.Item(26,6).Value = 14.2
$newws.Cells.Item(26,7).Value = 1
$newws.Cells.Item(26,8).Value = 0
$newws.Cells.Item(27,6).Value = 14.25
$newws.Cells.Item(27,7).Value = 1
$newws.Cells.Item(27,8).Value = 0
$newws.Cells.Item(28,6).Value = 14.3
$newws.Cells.Item(28,7).Value = 1
$newws.Cells.Item(28,8).Value = 0
$newws.Cells.Item(29,6).Value = 14.35
$newws.Cells.Item(29,7).Value = 1
$newws.Cells.Item(29,8).Value = 0
$newws.Cells.Item(30,6).Value = 14.4
$newws.Cells.Item(30,7).Value = 1
$newws.Cells.Item(30,8).Value = 0
$newws.Cells.Item(31,6).Value = 14.45
$newws.Cells.Item(31,7).Value = [double]"1.11e-100"
$newws.Cells.Item(31,8).Value = 21.26523355
$newws.Cells.Item(32,6).Value = 14.5
$newws.Cells.Item(32,7).Value = 1
$newws.Cells.Item(32,8).Value = 0
$newws.Cells.Item(33,6).Value = 14.55
$newws.Cells.Item(33,7).Value = 1
$newws.Cells.Item(33,8).Value = 0
$newws.Cells.Item(34,6).Value = 14.6
$newws.Cells.Item(34,7).Value = 1
$newws.Cells.Item(34,8).Value = 0
$newws.Cells.Item(35,6).Value = 14.65
$newws.Cells.Item(35,7).Value = 1
$newws.Cells.Item(35,8).Value = 0
$newws.Cells.Item(36,6).Value = 14.7
$newws.Cells.Item(36,7).Value = [double]"1.88e-108"
$newws.Cells.Item(36,8).Value = 22.08756955
$newws.Cells.Item(37,6).Value = 14.75
$newws.Cells.Item(37,7).Value = [double]"1.76e-108"
$newws.Cells.Item(37,8).Value = 22.08559389
$newws.Cells.Item(38,6).Value = 14.8
$newws.Cells.Item(38,7).Value = 1
$newws.Cells.Item(38,8).Value = 0
$newws.Cells.Item(39,6).Value = 14.85
$newws.Cells.Item(39,7).Value = 1
$newws.Cells.Item(39,8).Value = 0
$newws.Cells.Item(40,6).Value = 14.9
$newws.Cells.Item(40,7).Value = [double]"6.22e-206"
$newws.Cells.Item(40,8).Value = 30.59274947
$newws.Cells.Item(41,6).Value = 14.95
$newws.Cells.Item(41,7).Value = 1
$newws.Cells.Item(41,8).Value = 0
$newws.Cells.Item(42,6).Value = 15
$newws.Cells.Item(42,7).Value = 1
$newws.Cells.Item(42,8).Value = 0
$newws.Cells.Item(43,6).Value = 15.05
$newws.Cells.Item(43,7).Value = 1
$newws.Cells.Item(43,8).Value = 0
$newws.Cells.Item(44,6).Value = 15.1
$newws.Cells.Item(44,7).Value = 1
$newws.Cells.Item(44,8).Value = 0
$newws.Cells.Item(45,6).Value = 15.15
$newws.Cells.Item(45,7).Value = 0
$newws.Cells.Item(45,8).Value = 39.31677659
$newws.Cells.Item(46,6).Value = 15.2
$newws.Cells.Item(46,7).Value = [double]"2.97e-43"
$newws.Cells.Item(46,8).Value = 13.72634712
$newws.Cells.Item(47,6).Value = 15.25
$newws.Cells.Item(47,7).Value = 1
$newws.Cells.Item(47,8).Value = 0
$newws.Cells.Item(48,6).Value = 15.3
$newws.Cells.Item(48,7).Value = 1
$newws.Cells.Item(48,8).Value = 0
$newws.Cells.Item(49,6).Value = 15.35
$newws.Cells.Item(49,7).Value = [double]"1.53e-164"
$newws.Cells.Item(49,8).Value = 27.31068815
$newws.Cells.Item(50,6).Value = 15.4
$newws.Cells.Item(50,7).Value = [double]"3.44e-154"
$newws.Cells.Item(50,8).Value = 26.42082826
$newws.Cells.Item(51,6).Value = 15.45
$newws.Cells.Item(51,7).Value = 0
$newws.Cells.Item(51,8).Value = 40.48602588
$newws.Cells.Item(52,6).Value = 15.5
$newws.Cells.Item(52,7).Value = 1
$newws.Cells.Item(52,8).Value = 0
$newws.Cells.Item(53,6).Value = 15.55
$newws.Cells.Item(53,7).Value = 0
$newws.Cells.Item(53,8).Value = 38.85373053
$newws.Cells.Item(54,6).Value = 15.6
$newws.Cells.Item(54,7).Value = 1
$newws.Cells.Item(54,8).Value = 0
$newws.Cells.Item(55,6).Value = 15.65
$newws.Cells.Item(55,7).Value = [double]"5.63e-236"
$newws.Cells.Item(55,8).Value = 32.77440327
$newws.Cells.Item(56,6).Value = 15.7
$newws.Cells.Item(56,7).Value = [double]"1.03e-209"
$newws.Cells.Item(56,8).Value = 30.87977849
$newws.Cells.Item(57,6).Value = 15.75
$newws.Cells.Item(57,7).Value = [double]"1.27e-256"
$newws.Cells.Item(57,8).Value = 34.19724565
$newws.Cells.Item(58,6).Value = 15.8
$newws.Cells.Item(58,7).Value = 1
$newws.Cells.Item(58,8).Value = 0
$newws.Cells.Item(59,6).Value = 15.85
$newws.Cells.Item(59,7).Value = 0
$newws.Cells.Item(59,8).Value = 40.68737904
$newws.Cells.Item(60,6).Value = 15.9
$newws.Cells.Item(60,7).Value = 1
$newws.Cells.Item(60,8).Value = 0
$newws.Cells.Item(61,6).Value = 15.95
$newws.Cells.Item(61,7).Value = 0
$newws.Cells.Item(61,8).Value = 50.27601239
$newws.Cells.Item(62,6).Value = 16
$newws.Cells.Item(62,7).Value = 0
$newws.Cells.Item(62,8).Value = 40.18406275
$newws.Cells.Item(63,6).Value = 16.05
$newws.Cells.Item(63,7).Value = [double]"1.17e-292"
$newws.Cells.Item(63,8).Value = 36.54052376
$newws.Cells.Item(64,6).Value = 16.1
$newws.Cells.Item(64,7).Value = [double]"5.25e-230"
$newws.Cells.Item(64,8).Value = 32.35830427
$newws.Cells.Item(65,6).Value = 16.15
$newws.Cells.Item(65,7).Value = 0
$newws.Cells.Item(65,8).Value = 46.14361566
$newws.Cells.Item(66,6).Value = 16.2
$newws.Cells.Item(66,7).Value = 0
$newws.Cells.Item(66,8).Value = 48.78406207
$newws.Cells.Item(67,6).Value = 16.25
$newws.Cells.Item(67,7).Value = 0
$newws.Cells.Item(67,8).Value = 48.32187754
$newws.Cells.Item(68,6).Value = 16.3
$newws.Cells.Item(68,7).Value = 0
$newws.Cells.Item(68,8).Value = 49.63102249
$newws.Cells.Item(69,6).Value = 16.35
$newws.Cells.Item(69,7).Value = 0
$newws.Cells.Item(69,8).Value = 49.46844143
$newws.Cells.Item(70,6).Value = 16.4
$newws.Cells.Item(70,7).Value = 0
$newws.Cells.Item(70,8).Value = 52.66054594
$newws.Cells.Item(71,6).Value = 16.45
$newws.Cells.Item(71,7).Value = 0
$newws.Cells.Item(71,8).Value = 48.46359064
$newws.Cells.Item(72,6).Value = 16.5
$newws.Cells.Item(72,7).Value = 0
$newws.Cells.Item(72,8).Value = 50.26610726
$newws.Cells.Item(73,6).Value = 16.55
$newws.Cells.Item(73,7).Value = 0
$newws.Cells.Item(73,8).Value = 49.30188037
$newws.Cells.Item(74,6).Value = 16.6
$newws.Cells.Item(74,7).Value = 0
$newws.Cells.Item(74,8).Value = 54.5868235
$newws.Cells.Item(75,6).Value = 16.65
$newws.Cells.Item(75,7).Value = 0
$newws.Cells.Item(75,8).Value = 56.19857457
$newws.Cells.Item(76,6).Value = 16.7
$newws.Cells.Item(76,7).Value = 0
$newws.Cells.Item(76,8).Value = 48.00405242
$newws.Cells.Item(77,6).Value = 16.75
$newws.Cells.Item(77,7).Value = 0
$newws.Cells.Item(77,8).Value = 57.65872898
$newws.Cells.Item(78,6).Value = 16.8
$newws.Cells.Item(78,7).Value = 0
$newws.Cells.Item(78,8).Value = 55.1335156
$newws.Cells.Item(79,6).Value = 16.85
$newws.Cells.Item(79,7).Value = 0
$newws.Cells.Item(79,8).Value = 54.75949826
$newws.Cells.Item(80,6).Value = 16.9
$newws.Cells.Item(80,7).Value = 0
$newws.Cells.Item(80,8).Value = 56.17374498
$newws.Cells.Item(81,6).Value = 16.95
$newws.Cells.Item(81,7).Value = 0
$newws.Cells.Item(81,8).Value = 61.16670593
$newws.Cells.Item(82,6).Value = 17
$newws.Cells.Item(82,7).Value = 0
$newws.Cells.Item(82,8).Value = 58.00233029
$newws.Cells.Item(83,6).Value = 17.05
$newws.Cells.Item(83,7).Value = 0
$newws.Cells.Item(83,8).Value = 51.67073046
$newws.Cells.Item(84,6).Value = 17.1
$newws.Cells.Item(84,7).Value = 0
$newws.Cells.Item(84,8).Value = 53.79423895
$newws.Cells.Item(85,6).Value = 17.15
$newws.Cells.Item(85,7).Value = 0
$newws.Cells.Item(85,8).Value = 62.20640968
$newws.Cells.Item(86,6).Value = 17.2
$newws.Cells.Item(86,7).Value = 0
$newws.Cells.Item(86,8).Value = 62.61363283
$newws.Cells.Item(87,6).Value = 17.25
$newws.Cells.Item(87,7).Value = 0
$newws.Cells.Item(87,8).Value = 63.69004851
$newws.Cells.Item(88,6).Value = 17.3
$newws.Cells.Item(88,7).Value = 0
$newws.Cells.Item(88,8).Value = 47.30045527
$newws.Cells.Item(89,6).Value = 17.35
$newws.Cells.Item(89,7).Value = 0
$newws.Cells.Item(89,8).Value = 58.13297972
$newws.Cells.Item(90,6).Value = 17.4
$newws.Cells.Item(90,7).Value = 0
$newws.Cells.Item(90,8).Value = 57.94121675
$newws.Cells.Item(91,6).Value = 17.45
$newws.Cells.Item(91,7).Value = 0
$newws.Cells.Item(91,8).Value = 58.1182361
$newws.Cells.Item(92,6).Value = 17.5
$newws.Cells.Item(92,7).Value = 0
$newws.Cells.Item(92,8).Value = 56.37006178
$newws.Cells.Item(93,6).Value = 17.55
$newws.Cells.Item(93,7).Value = 0
$newws.Cells.Item(93,8).Value = 65.74543345
$newws.Cells.Item(94,6).Value = 17.6
$newws.Cells.Item(94,7).Value = 0
$newws.Cells.Item(94,8).Value = 64.37381271
$newws.Cells.Item(95,6).Value = 17.65
$newws.Cells.Item(95,7).Value = 0
$newws.Cells.Item(95,8).Value = 63.07811396
$newws.Cells.Item(96,6).Value = 17.7
$newws.Cells.Item(96,7).Value = 0
$newws.Cells.Item(96,8).Value = 63.40377858
$newws.Cells.Item(97,6).Value = 17.75
$newws.Cells.Item(97,7).Value = 0
$newws.Cells.Item(97,8).Value = 69.78771301
$newws.Cells.Item(98,6).Value = 17.8
$newws.Cells.Item(98,7).Value = 0
$newws.Cells.Item(98,8).Value = 59.47437001
$newws.Cells.Item(99,6).Value = 17.85
$newws.Cells.Item(99,7).Value = 0
$newws.Cells.Item(99,8).Value = 57.37299158
$newws.Cells.Item(100,6).Value = 17.9
$newws.Cells.Item(100,7).Value = 0
$newws.Cells.Item(100,8).Value = 67.55811829
$newws.Cells.Item(101,6).Value = 17.95
$newws.Cells.Item(101,7).Value = 0
$newws.Cells.Item(101,8).Value = 65.24354974
$newws.Cells.Item(102,6).Value = 18
$newws.Cells.Item(102,7).Value = 0
$newws.Cells.Item(102,8).Value = 64.30596057

$newws.Cells.Item(31,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(36,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(37,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(40,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(46,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(49,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(50,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(55,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(56,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(57,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(63,7).NumberFormat = "0.00E+00"
$newws.Cells.Item(64,7).NumberFormat = "0.00E+00"

$newws.Activate()
$newws.Range("E8").Select()
